$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shared string "h_t" (girl's hat/hit) is used by cells B37 and B41.
# Update it in both places to "dr_ss" (dress). Changing the shared string
# text itself (rather than re-pointing to a different entry) reproduces
# the sharedStrings.xml diff, since "h_t" -> "dr_ss" is edited in place.
$ws.Range("B37").Value = "dr_ss"
$ws.Range("B41").Value = "dr_ss"

# The corresponding correct-answer letters change from "a" to "e".
$ws.Range("C37").Value = "e"
$ws.Range("C41").Value = "e"

# Update the view state: scroll so row 28 is the top-left visible row,
# and select cell B41 as the active cell.
$ws.Range("B41").Select()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
